# "Register" flow: a new user (matan123) successfully registered and was
# appended as a new row under the existing Users table, mirroring what the
# app's Register button now writes to Users.xlsx before returning the user
# to the login form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New user row (row 3): Username / Full name / Password -----------------
$ws.Cells.Item(3, 1).Value = "matan123"
$ws.Cells.Item(3, 2).Value = "matan123@$"

# The password "315783522" looks like a pure number, so Excel would normally
# coerce it to a numeric cell. Force it to be stored as text (matching the
# other, textual, "Password" entries) via a text-formatted helper cell that
# we paste the *value* of (so the destination keeps the sheet's normal
# styling) and then discard.
$helper = $ws.Cells.Item(100, 100)
$helper.NumberFormat = "@"
$helper.Value = "315783522"
$helper.Copy()
$ws.Cells.Item(3, 3).PasteSpecial(-4163)  # xlPasteValues
$helper.Clear()

# --- Trailing blank rows (4 and 5), kept in the sheet's used range ---------
foreach ($r in 4, 5) {
    foreach ($c in 1, 2, 3) {
        $cell = $ws.Cells.Item($r, $c)
        # Touch (no-op) formatting so the cell is materialized in the sheet
        # without a value and without changing its (default) style.
        $cell.Font.Size = 11
    }
}

# --- Page setup / print settings -------------------------------------------
$ws.PageSetup.PaperSize = 1
$ws.PageSetup.Zoom = 100
$ws.PageSetup.CenterHeader = $ws.PageSetup.CenterHeader
